# Terra Machina Reference Guide - "Trying to fix name generator"
$wb = $excel.ActiveWorkbook

# --- Rename sheet 3 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Sample_Custom_Moves"

# --- Sheet1 (Items): new rows 15-17 (Potion+, Tincture+, Concoction+) ---
$ws1.Range("A15").Value = "Potion+"
$ws1.Range("B15").Value = "Heals All HP"
$ws1.Range("C15").Value = 25

$ws1.Range("A16").Value = "Tincture+"
$ws1.Range("B16").Value = "Heals all HP and SP"
$ws1.Range("C16").Value = 70

$ws1.Range("A17").Value = "Concoction+"
$ws1.Range("B17").Value = "Heals all SP"
$ws1.Range("C17").Value = 25

# --- Sheet1 (Items): row 7 "Tonic" -> "Brew" ---
$ws1.Range("A7").Value = "Brew"

# --- Sheet1 (Items): new row 18 (Lazarus Tonic) ---
$ws1.Range("A18").Value = "Lazarus Tonic"
$ws1.Range("B18").Value = "Revives a fainted comrade in battle"
$ws1.Range("C18").Value = 100

# --- Sheet2 (Bestiary): new row 6 "Gunner Soldier" ---
$ws2.Range("A6").Value = "Gunner Soldier"
$ws2.Range("B6").Value = "A soldier with a gun"
$ws2.Range("C6").Value = 7
$ws2.Range("D6").Value = 4
$ws2.Range("F6").Value = "Gunshot: 1d10 damage to an enemy`nGun Spray: 1d4 damage to all enemies"

# --- Sheet3 (Sample_Custom_Moves): new row 3 "Prosecutor" ---
$ws3.Range("A3").Value = "Prosecutor"
$ws3.Range("B3").Value = "'+1 to Diplomacy when you try to convince others that someone is lying (whether it's true or not)"

# --- Sheet2 (Bestiary): Armored Rhino weakness format change (last) ---
$ws2.Range("E3").Value = "Fire, Lightning"

# --- Selections: update non-active sheets first, then select on the
#     active sheet (Items) last so it keeps focus / tabSelected ---
$ws2.Range("A8").Select()
$ws3.Range("A4").Select()
$ws1.Range("G13").Select()
